$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value2 = 13.73601491863699
$ws.Range("C2").Value2 = 10.66670083838895
$ws.Range("D2").Value2 = 5.977444122007769
$ws.Range("E2").Value2 = 12.47962550562971
$ws.Range("G2").Value2 = 27.65683162255755
$ws.Range("H2").Value2 = 13.89745610560035
$ws.Range("I2").Value2 = 21.46468702686644
$ws.Range("L2").Value2 = 9.998992387164492
$ws.Range("M2").Value2 = 14.12300526268736
$ws.Range("N2").Value2 = 17.64873399393152
$ws.Range("O2").Value2 = 21.02916406240776

$ws.Range("B3").Value2 = 13.25997561723036
$ws.Range("C3").Value2 = 10.46598940404726
$ws.Range("D3").Value2 = 5.857140931170137
$ws.Range("E3").Value2 = 12.51923234144846
$ws.Range("G3").Value2 = 27.59755330861018
$ws.Range("H3").Value2 = 13.93447248047569
$ws.Range("I3").Value2 = 21.5555621622159
$ws.Range("L3").Value2 = 10.00629420462095
$ws.Range("M3").Value2 = 14.02789378208334
$ws.Range("N3").Value2 = 17.69608465666053
$ws.Range("O3").Value2 = 21.06786359910546

$ws.Range("B4").Value2 = 12.96030291970072
$ws.Range("C4").Value2 = 10.33966492999481
$ws.Range("D4").Value2 = 5.783770203523147
$ws.Range("E4").Value2 = 12.54489812231173
$ws.Range("G4").Value2 = 27.5711067009284
$ws.Range("H4").Value2 = 13.95973375719857
$ws.Range("I4").Value2 = 21.61581865088586
$ws.Range("L4").Value2 = 10.01217586479951
$ws.Range("M4").Value2 = 13.97113265224161
$ws.Range("N4").Value2 = 17.72688411377214
$ws.Range("O4").Value2 = 21.09682146171697

$ws.Range("B5").Value2 = 12.83651223581193
$ws.Range("C5").Value2 = 10.28745213611728
$ws.Range("D5").Value2 = 5.754042536225199
$ws.Range("E5").Value2 = 12.55569667508282
$ws.Range("G5").Value2 = 27.5628384596183
$ws.Range("H5").Value2 = 13.97066422373562
$ws.Range("I5").Value2 = 21.64149329060131
$ws.Range("L5").Value2 = 10.01492484911594
$ws.Range("M5").Value2 = 13.94843253508186
$ws.Range("N5").Value2 = 17.73987001921327
$ws.Range("O5").Value2 = 21.10992577034045

$ws.Range("B6").Value2 = 12.81586166113606
$ws.Range("C6").Value2 = 10.27873913011507
$ws.Range("D6").Value2 = 5.749118015476327
$ws.Range("E6").Value2 = 12.55751029798247
$ws.Range("G6").Value2 = 27.56161719526983
$ws.Range("H6").Value2 = 13.97251762924455
$ws.Range("I6").Value2 = 21.64582412236453
$ws.Range("L6").Value2 = 10.01540259986734
$ws.Range("M6").Value2 = 13.94468973267796
$ws.Range("N6").Value2 = 17.74205261023456
$ws.Range("O6").Value2 = 21.1121803814882

$ws.Range("B7").Value2 = 12.95863994684603
$ws.Range("C7").Value2 = 10.33896368775091
$ws.Range("D7").Value2 = 5.783368529179629
$ws.Range("E7").Value2 = 12.54504237949619
$ws.Range("G7").Value2 = 27.57098502725337
$ws.Range("H7").Value2 = 13.95987859390415
$ws.Range("I7").Value2 = 21.61616037678285
$ws.Range("L7").Value2 = 10.01221151201686
$ws.Range("M7").Value2 = 13.97082474322691
$ws.Range("N7").Value2 = 17.72705748427428
$ws.Range("O7").Value2 = 21.09699291632008

$ws.Range("B8").Value2 = 13.5735144208981
$ws.Range("C8").Value2 = 10.59816062611575
$ws.Range("D8").Value2 = 5.935888536381806
$ws.Range("E8").Value2 = 12.49300286138005
$ws.Range("G8").Value2 = 27.63433231491215
$ws.Range("H8").Value2 = 13.90969322347473
$ws.Range("I8").Value2 = 21.49509441061269
$ws.Range("L8").Value2 = 10.00122017913682
$ws.Range("M8").Value2 = 14.08988232077053
$ws.Range("N8").Value2 = 17.66470280141353
$ws.Range("O8").Value2 = 21.04142758232453

$ws.Range("B9").Value2 = 14.71344706668401
$ws.Range("C9").Value2 = 11.08022090162224
$ws.Range("D9").Value2 = 6.236897990123045
$ws.Range("E9").Value2 = 12.40160408543471
$ws.Range("G9").Value2 = 27.83705591997076
$ws.Range("H9").Value2 = 13.83140723149882
$ws.Range("I9").Value2 = 21.29313286167088
$ws.Range("L9").Value2 = 9.990735301509906
$ws.Range("M9").Value2 = 14.33550007849351
$ws.Range("N9").Value2 = 17.55608094381513
$ws.Range("O9").Value2 = 20.9738056956536

$ws.Range("B10").Value2 = 15.50227118787785
$ws.Range("C10").Value2 = 11.41623060142122
$ws.Range("D10").Value2 = 6.456604875523861
$ws.Range("E10").Value2 = 12.34089368905028
$ws.Range("G10").Value2 = 28.03302525697008
$ws.Range("H10").Value2 = 13.78619758394184
$ws.Range("I10").Value2 = 21.16645762823542
$ws.Range("L10").Value2 = 9.98974216182431
$ws.Range("M10").Value2 = 14.52219042227895
$ws.Range("N10").Value2 = 17.4845470207644
$ws.Range("O10").Value2 = 20.94945508406273

$ws.Range("B11").Value2 = 15.84906366865787
$ws.Range("C11").Value2 = 11.56472890148346
$ws.Range("D11").Value2 = 6.555721499343282
$ws.Range("E11").Value2 = 12.31466226952868
$ws.Range("G11").Value2 = 28.13214460899416
$ws.Range("H11").Value2 = 13.7683097263319
$ws.Range("I11").Value2 = 21.11356381005023
$ws.Range("L11").Value2 = 9.990737725780589
$ws.Range("M11").Value2 = 14.60821743240595
$ws.Range("N11").Value2 = 17.45378902296575
$ws.Range("O11").Value2 = 20.94389560347884

$ws.Range("B12").Value2 = 15.97854740328699
$ws.Range("C12").Value2 = 11.62030308300346
$ws.Range("D12").Value2 = 6.593094657066405
$ws.Range("E12").Value2 = 12.30492759865453
$ws.Range("G12").Value2 = 28.17108715402028
$ws.Range("H12").Value2 = 13.76192164784656
$ws.Range("I12").Value2 = 21.09421647902836
$ws.Range("L12").Value2 = 9.991321867151779
$ws.Range("M12").Value2 = 14.64092983399587
$ws.Range("N12").Value2 = 17.44239737127993
$ws.Range("O12").Value2 = 20.94258439618514

$ws.Range("B13").Value2 = 15.95074411568284
$ws.Range("C13").Value2 = 11.60836399211865
$ws.Range("D13").Value2 = 6.585053428613313
$ws.Range("E13").Value2 = 12.30701531274828
$ws.Range("G13").Value2 = 28.16263799147679
$ws.Range("H13").Value2 = 13.76328027447082
$ws.Range("I13").Value2 = 21.09835288999104
$ws.Range("L13").Value2 = 9.99118686572219
$ws.Range("M13").Value2 = 14.633878949957
$ws.Range("N13").Value2 = 17.44483940403563
$ws.Range("O13").Value2 = 20.94283146866831

$ws.Range("B14").Value2 = 15.85975379502254
$ws.Range("C14").Value2 = 11.56931440817779
$ws.Range("D14").Value2 = 6.558799647409029
$ws.Range("E14").Value2 = 12.31385741614199
$ws.Range("G14").Value2 = 28.13532037119311
$ws.Range("H14").Value2 = 13.7677764431375
$ws.Range("I14").Value2 = 21.11195840525431
$ws.Range("L14").Value2 = 9.990781637938113
$ws.Range("M14").Value2 = 14.61090610786661
$ws.Range("N14").Value2 = 17.45284670344473
$ws.Range("O14").Value2 = 20.94377181419127

$ws.Range("B15").Value2 = 15.8037771569861
$ws.Range("C15").Value2 = 11.54530866662448
$ws.Range("D15").Value2 = 6.542696416418262
$ws.Range("E15").Value2 = 12.31807424439591
$ws.Range("G15").Value2 = 28.11877011427123
$ws.Range("H15").Value2 = 13.77058071600617
$ws.Range("I15").Value2 = 21.12038110683234
$ws.Range("L15").Value2 = 9.99056036964082
$ws.Range("M15").Value2 = 14.59685160119167
$ws.Range("N15").Value2 = 17.45778468990452
$ws.Range("O15").Value2 = 20.94445122030252

$ws.Range("B16").Value2 = 15.47935563843952
$ws.Range("C16").Value2 = 11.4064354192036
$ws.Range("D16").Value2 = 6.450107135190272
$ws.Range("E16").Value2 = 12.34263580360207
$ws.Range("G16").Value2 = 28.02674610119033
$ws.Range("H16").Value2 = 13.78742053802926
$ws.Range("I16").Value2 = 21.17000974568482
$ws.Range("L16").Value2 = 9.989706141611665
$ws.Range("M16").Value2 = 14.51658858836146
$ws.Range("N16").Value2 = 17.48659295113418
$ws.Range("O16").Value2 = 20.9499295037415

$ws.Range("B17").Value2 = 15.27716776367006
$ws.Range("C17").Value2 = 11.3201019529318
$ws.Range("D17").Value2 = 6.393064397676985
$ws.Range("E17").Value2 = 12.3580580136221
$ws.Range("G17").Value2 = 27.97282902399355
$ws.Range("H17").Value2 = 13.79843758868063
$ws.Range("I17").Value2 = 21.20166858059095
$ws.Range("L17").Value2 = 9.989552128681417
$ws.Range("M17").Value2 = 14.46761581933734
$ws.Range("N17").Value2 = 17.50472210237986
$ws.Range("O17").Value2 = 20.95470407355676

$ws.Range("B18").Value2 = 15.15974616364716
$ws.Range("C18").Value2 = 11.2700372093828
$ws.Range("D18").Value2 = 6.360179187158925
$ws.Range("E18").Value2 = 12.3670589640559
$ws.Range("G18").Value2 = 27.94275768667187
$ws.Range("H18").Value2 = 13.80502636534276
$ws.Range("I18").Value2 = 21.22032312247912
$ws.Range("L18").Value2 = 9.989599840194067
$ws.Range("M18").Value2 = 14.43955292936968
$ws.Range("N18").Value2 = 17.5153173878324
$ws.Range("O18").Value2 = 20.95796963200592

$ws.Range("B19").Value2 = 15.11979887663284
$ws.Range("C19").Value2 = 11.25301710452296
$ws.Range("D19").Value2 = 6.349033096938827
$ws.Range("E19").Value2 = 12.37012896712449
$ws.Range("G19").Value2 = 27.9327382885803
$ws.Range("H19").Value2 = 13.80730048605237
$ws.Range("I19").Value2 = 21.22671563418764
$ws.Range("L19").Value2 = 9.989639431951952
$ws.Range("M19").Value2 = 14.43007002352089
$ws.Range("N19").Value2 = 17.51893362352448
$ws.Range("O19").Value2 = 20.95916446624789

$ws.Range("B20").Value2 = 15.29880861600074
$ws.Range("C20").Value2 = 11.32933477038224
$ws.Range("D20").Value2 = 6.399144836769604
$ws.Range("E20").Value2 = 12.35640279165037
$ws.Range("G20").Value2 = 27.97847144440949
$ws.Range("H20").Value2 = 13.79723871339116
$ws.Range("I20").Value2 = 21.19825235044139
$ws.Range("L20").Value2 = 9.989554424380973
$ws.Range("M20").Value2 = 14.47281835414375
$ws.Range("N20").Value2 = 17.50277485360188
$ws.Range("O20").Value2 = 20.95414205539116

$ws.Range("B21").Value2 = 15.88653055408107
$ws.Range("C21").Value2 = 11.58080235032315
$ws.Range("D21").Value2 = 6.566515680525859
$ws.Range("E21").Value2 = 12.31184234025558
$ws.Range("G21").Value2 = 28.14330622019209
$ws.Range("H21").Value2 = 13.76644533863373
$ws.Range("I21").Value2 = 21.1079435989052
$ws.Range("L21").Value2 = 9.990895049639159
$ws.Range("M21").Value2 = 14.61765028118721
$ws.Range("N21").Value2 = 17.45048783026104
$ws.Range("O21").Value2 = 20.94347405950227

$ws.Range("B22").Value2 = 16.25988095788628
$ws.Range("C22").Value2 = 11.74129669991432
$ws.Range("D22").Value2 = 6.67495137678344
$ws.Range("E22").Value2 = 12.2838767682091
$ws.Range("G22").Value2 = 28.25922934235178
$ws.Range("H22").Value2 = 13.74856823858499
$ws.Range("I22").Value2 = 21.05290016705071
$ws.Range("L22").Value2 = 9.992978094671161
$ws.Range("M22").Value2 = 14.71308762242976
$ws.Range("N22").Value2 = 17.41780555334756
$ws.Range("O22").Value2 = 20.94113012874129

$ws.Range("B23").Value2 = 16.06163324469687
$ws.Range("C23").Value2 = 11.65600077327597
$ws.Range("D23").Value2 = 6.61717711704584
$ws.Range("E23").Value2 = 12.29869686434377
$ws.Range("G23").Value2 = 28.19661849220446
$ws.Range("H23").Value2 = 13.75790371786843
$ws.Range("I23").Value2 = 21.0819131683224
$ws.Range("L23").Value2 = 9.991756253504057
$ws.Range("M23").Value2 = 14.66208686638851
$ws.Range("N23").Value2 = 17.43511255495363
$ws.Range("O23").Value2 = 20.9419575834522

$ws.Range("B24").Value2 = 15.28902845727254
$ws.Range("C24").Value2 = 11.32516195199918
$ws.Range("D24").Value2 = 6.396396149270762
$ws.Range("E24").Value2 = 12.35715069844929
$ws.Range("G24").Value2 = 27.97591761883275
$ws.Range("H24").Value2 = 13.7977799308576
$ws.Range("I24").Value2 = 21.19979541574752
$ws.Range("L24").Value2 = 9.989552962002056
$ws.Range("M24").Value2 = 14.47046599922705
$ws.Range("N24").Value2 = 17.50365466709275
$ws.Range("O24").Value2 = 20.95439452222737

$ws.Range("B25").Value2 = 14.41305828772298
$ws.Range("C25").Value2 = 10.95285076899542
$ws.Range("D25").Value2 = 6.155537603762181
$ws.Range("E25").Value2 = 12.42519512820423
$ws.Range("G25").Value2 = 27.77388237179303
$ws.Range("H25").Value2 = 13.85042685484618
$ws.Range("I25").Value2 = 21.34396398243099
$ws.Range("L25").Value2 = 9.992390302774156
$ws.Range("M25").Value2 = 14.26787263355308
$ws.Range("N25").Value2 = 17.58400979669399
$ws.Range("O25").Value2 = 20.98765820176654

